# issue #5: add legislator_id, name, date into dataframe
#
# Adds three new columns (date, legislator_name, legislator_id) to the
# "股票" (stock) sheet, filled with the report date, the legislator's
# name, and their id for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

$legislatorName = "王育敏"
$legislatorId = 1728
$reportDate = "2013-11-12"

# Header row (row 1) - same bold/centered/bordered style as existing headers (copy from G1)
$ws.Range("H1").Value = "date"
$ws.Range("I1").Value = "legislator_name"
$ws.Range("J1").Value = "legislator_id"
$ws.Range("H1:J1").Style = $ws.Range("G1").Style

# Data rows - same plain style as the rest of the data cells (copy from G2/G3)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 8).Value = $reportDate
    $ws.Cells.Item($r, 9).Value = $legislatorName
    $ws.Cells.Item($r, 10).Value = $legislatorId
    $ws.Range($ws.Cells.Item($r, 8), $ws.Cells.Item($r, 10)).Style = $ws.Cells.Item($r, 7).Style
}
